# Update cryptocurrency price/volume snapshot cells (columns D and E)
# per the latest GitHub Actions scrape. Price cells that are purely
# numeric-looking text (a single decimal point) are forced to the "@"
# (Text) number format before assignment so Excel keeps them as literal
# strings instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.479.72"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.837.20"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.98"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5391"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3012"
$ws.Range("E8").Value = "  -7.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06913"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.63"
$ws.Range("E10").Value = "  -6.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7381"
$ws.Range("D12").Value = "1.846.66"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07146"
$ws.Range("E13").Value = "  -8.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.25"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.993"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.85"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007897"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "26.509.94"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "2.072.14"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.585"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.988"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.209"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.74"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.178"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.01"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.01"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.249"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08853"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04843"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.933"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7307"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.092"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4713"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9060"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.10"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.902"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.406"
$ws.Range("E45").Value = "  -3.52%  "
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.014"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4080"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.84"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8945"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05765"
$ws.Range("E51").Value = "  -2.07%  "
